# Update gh-pages output — applies the data refresh captured in the commit.
#
# Summary of the change:
#  - Sheet "展览" (exhibitions): several "want-to-go" counts (column F) bumped.
#  - Sheet "演出" (performances): one "want-to-go" count (column F) bumped.
#  - Sheet "本地生活" (local life): one "want-to-go" count (column F) bumped.
#  - Sheet "全部类型" (all types, a date-sorted merge of the other three
#    sheets): a brand-new performance row (BELOVED concert on 2024-09-16)
#    was inserted at row 10, pushing the following rows down by one; a
#    pre-existing duplicate row (the "GZDP10" doll exhibition, which used
#    to appear twice back-to-back) was removed, which keeps every row
#    below the insertion point at its original row number. On top of
#    that shuffle, several of those same rows also got their column-F
#    counts bumped, and a few untouched rows further down (>=22) only had
#    their column-F counts bumped.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "展览" — simple column F bumps
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F9").Value  = 2136
$wsExpo.Range("F10").Value = 614
$wsExpo.Range("F12").Value = 112
$wsExpo.Range("F13").Value = 1039
$wsExpo.Range("F15").Value = 2176
$wsExpo.Range("F16").Value = 632
$wsExpo.Range("F17").Value = 11725
$wsExpo.Range("F18").Value = 1184
$wsExpo.Range("F21").Value = 6
$wsExpo.Range("F22").Value = 135
$wsExpo.Range("F27").Value = 9

# ---------------------------------------------------------------------
# 2) Sheet "演出" — simple column F bump
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F10").Value = 14

# ---------------------------------------------------------------------
# 3) Sheet "本地生活" — simple column F bump
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5678

# ---------------------------------------------------------------------
# 4) Sheet "全部类型" — row 2 (本地生活 dup) F bump, then the row
#    insertion/shuffle for rows 10-21, then the remaining untouched-row F bumps.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

# Same 本地生活 count bump, duplicated onto this merged sheet (row 3 here).
$wsAll.Range("F3").Value = 5678

# Rows 10-21 get fully rewritten in place: a new row is inserted at 10
# (shifting old rows 10-20 down to 11-21) while the old duplicate that used
# to sit at row 21 is dropped, so every row from 22 onward keeps its
# original row number. Column A (running index) and the header row are
# untouched throughout.

# Row 10 (NEW): 广州·中秋《BELOVED·喜欢你》华语三十年经典歌曲演唱会
$wsAll.Range("B10").Value = "2024-09-16"
$wsAll.Range("C10").Value = "广州·中秋《BELOVED·喜欢你》华语三十年经典歌曲演唱会"
$wsAll.Range("D10").Value = "香雪大道西3号 广州科学城会议中心"
$wsAll.Range("E10").Value = "2024.09.16 19:30-09.16 21:10"
$wsAll.Range("F10").Value = 2
$wsAll.Range("G10").Value = 136
$wsAll.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=91806"
$wsAll.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202409/ZnA9tjLO1725328334894.png"

# Row 11: 广州·星光次元动漫嘉年华 (was row 10)
$wsAll.Range("B11").Value = "2024-09-16"
$wsAll.Range("C11").Value = "广州·星光次元动漫嘉年华"
$wsAll.Range("D11").Value = "奥体南路12号 优托邦购物中心(奥体中心店)"
$wsAll.Range("E11").Value = "2024.09.16 10:00-09.16 17:00"
$wsAll.Range("F11").Value = 457
$wsAll.Range("G11").Value = 39.9
$wsAll.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=88949"
$wsAll.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202407/fmuLNRYC1719799438032.jpeg"

# Row 12: 广州·经典重现·电影配乐交响音乐会 (was row 11)
$wsAll.Range("B12").Value = "2024-09-16"
$wsAll.Range("C12").Value = "广州·经典重现·电影配乐交响音乐会"
$wsAll.Range("D12").Value = "晴波路33号 广州星海音乐厅"
$wsAll.Range("E12").Value = "2024.09.16 20:00-09.16 21:40"
$wsAll.Range("F12").Value = 15
$wsAll.Range("G12").Value = 126
$wsAll.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=89743"
$wsAll.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202407/vFb5Ojv91721620888814.jpeg"

# Row 13: 广州·首届China Comic次元动漫游戏嘉年华 (was row 12, F bumped 2130->2136)
$wsAll.Range("B13").Value = "2024-09-16"
$wsAll.Range("C13").Value = "广州·首届China Comic次元动漫游戏嘉年华"
$wsAll.Range("D13").Value = "广汕二路149号 飞梦篮球公园(高塘石馆)"
$wsAll.Range("E13").Value = "2024.09.16 10:00-09.17 17:00"
$wsAll.Range("F13").Value = 2136
$wsAll.Range("G13").Value = 60
$wsAll.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=91595"
$wsAll.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202408/WNDf0pdH1725000637494.jpeg"

# Row 14: 广州·樱漫潮玩动漫嘉年华 (was row 13, F bumped 611->614)
$wsAll.Range("B14").Value = "2024-09-17"
$wsAll.Range("C14").Value = "广州·樱漫潮玩动漫嘉年华"
$wsAll.Range("D14").Value = "开创大道2666号 宝能国际体育演艺中心"
$wsAll.Range("E14").Value = "2024.09.17 10:00-09.17 17:00"
$wsAll.Range("F14").Value = 614
$wsAll.Range("G14").Value = 39.9
$wsAll.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=88936"
$wsAll.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202407/v2XniRgx1719799991150.jpeg"

# Row 15: 广州·VOCALOID术力口同人only-2.0 (was row 14)
$wsAll.Range("B15").Value = "2024-09-21"
$wsAll.Range("C15").Value = "广州·VOCALOID术力口同人only-2.0"
$wsAll.Range("D15").Value = "华观路1932号 智慧城广百广场"
$wsAll.Range("E15").Value = "2024.09.21 10:00-09.22 17:00"
$wsAll.Range("F15").Value = 276
$wsAll.Range("G15").Value = 50
$wsAll.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=90942"
$wsAll.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202408/HieGJX2e1723786754177.jpeg"

# Row 16: 广州·《梁祝·卡农》中外经典名曲精选音乐会 (was row 15)
$wsAll.Range("B16").Value = "2024-09-21"
$wsAll.Range("C16").Value = "广州·《梁祝·卡农》中外经典名曲精选音乐会"
$wsAll.Range("D16").Value = "凤凰北路65号 融创广州大剧院"
$wsAll.Range("E16").Value = "2024.09.21 19:30-09.21 21:00"
$wsAll.Range("F16").Value = 13
$wsAll.Range("G16").Value = 38
$wsAll.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=89959"
$wsAll.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202407/MVreW47V1721899542180.jpeg"

# Row 17: 广州·守望先锋同人ONLY (was row 16, F bumped 111->112)
$wsAll.Range("B17").Value = "2024-09-21"
$wsAll.Range("C17").Value = "广州·守望先锋同人ONLY"
$wsAll.Range("D17").Value = "西环路1号 广州岭南会展中心"
$wsAll.Range("E17").Value = "2024.09.21 09:00-09.21 16:30"
$wsAll.Range("F17").Value = 112
$wsAll.Range("G17").Value = 100
$wsAll.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=91569"
$wsAll.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202408/Dgan6x331724939875504.jpeg"

# Row 18: 广州·VGL中国巡演 2024 VIDEO GAMES LIVE 魔兽世界音乐会 (was row 17)
$wsAll.Range("B18").Value = "2024-09-22"
$wsAll.Range("C18").Value = "广州·VGL中国巡演 2024 VIDEO GAMES LIVE 魔兽世界音乐会"
$wsAll.Range("D18").Value = "东风中路299号 广州中山纪念堂"
$wsAll.Range("E18").Value = "2024.09.22 19:30-09.22 21:10"
$wsAll.Range("F18").Value = 22
$wsAll.Range("G18").Value = 180
$wsAll.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=88919"
$wsAll.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202407/Bn0OQ6ef1720428966541.jpeg"

# Row 19: 广州·wio流金序曲乙女同人ONLY展 (was row 18, F bumped 1034->1039)
$wsAll.Range("B19").Value = "2024-09-22"
$wsAll.Range("C19").Value = "广州·wio流金序曲乙女同人ONLY展"
$wsAll.Range("D19").Value = "黄边三横路一街1号 设计殿堂"
$wsAll.Range("E19").Value = "2024.09.22 10:00-09.22 21:00"
$wsAll.Range("F19").Value = 1039
$wsAll.Range("G19").Value = 69.90000000000001
$wsAll.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=89569"
$wsAll.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202407/edi1S6tl1721231532813.png"

# Row 20: 广州·傻子与白痴乐队「姿态」2024全国巡演 (was row 19)
$wsAll.Range("B20").Value = "2024-09-22"
$wsAll.Range("C20").Value = "广州·傻子与白痴乐队「姿态」2024全国巡演"
$wsAll.Range("D20").Value = "龙凤街道革新路124号太古仓码头5号仓 广州太空间Live House"
$wsAll.Range("E20").Value = "2024.09.22 20:00-09.22 22:00"
$wsAll.Range("F20").Value = 0
$wsAll.Range("G20").Value = 299
$wsAll.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=91966"
$wsAll.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202409/E0wfaDqh1725504323717.jpeg"

# Row 21: 广州·第10届娃展沙龙-顽哇展(GZDP10) (was row 20; the old row-21
# duplicate of this same event is the one that got dropped)
$wsAll.Range("B21").Value = "2024-09-22"
$wsAll.Range("C21").Value = "广州·第10届娃展沙龙-顽哇展(GZDP10)"
$wsAll.Range("D21").Value = "洛浦街厦滘西环路1号 岭南国际电子商务会展中心"
$wsAll.Range("E21").Value = "2024.09.22 11:00-09.22 17:00"
$wsAll.Range("F21").Value = 172
$wsAll.Range("G21").Value = 55
$wsAll.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=90741"
$wsAll.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202408/R8lHU3dM1723534719530.png"

# Rows 22+ keep their original row numbers and content; only a handful of
# them get their column F "want-to-go" counts bumped.
$wsAll.Range("F23").Value = 14
$wsAll.Range("F24").Value = 2176
$wsAll.Range("F25").Value = 632
$wsAll.Range("F26").Value = 11726
$wsAll.Range("F29").Value = 1184
$wsAll.Range("F32").Value = 6
$wsAll.Range("F33").Value = 135
$wsAll.Range("F49").Value = 9
